$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "{0: sala nr 9 | Natalia Szymańska | Geografia}"
$ws.Range("C2").Value = "{0: sala nr 7 | Paweł Lewandowski | Matematyka}"
$ws.Range("D2").Value = "{0: sala nr 1 | Paweł Lewandowski | Matematyka}"
$ws.Range("E2").Value = "{}"
$ws.Range("F2").Value = "{}"

$ws.Range("B3").Value = "{0: sala nr 2 | Lena Kowalska | Język angielski}"
$ws.Range("C3").Value = "{0: sala nr 7 | Paweł Lewandowski | Matematyka}"
$ws.Range("D3").Value = "{0: sala nr 5 | Dominik Kaczor | Informatyka}"
$ws.Range("E3").Value = "{}"
$ws.Range("F3").Value = "{}"

$ws.Range("B4").Value = "{0: sala nr 5 | Katarzyna Mazur | Fizyka}"
$ws.Range("C4").Value = "{0: sala nr 6 | Katarzyna Mazur | Fizyka}"
$ws.Range("D4").Value = "{}"
$ws.Range("E4").Value = "{}"
$ws.Range("F4").Value = "{}"

$ws.Range("B5").Value = "{0: sala nr 5 | Piotr Wójcik | Biologia}"
$ws.Range("C5").Value = "{0: sala nr 10 | Jan Nowak | Język polski}"
$ws.Range("D5").Value = "{0: sala nr 4 | Lena Kowalska | Język angielski}"
$ws.Range("E5").Value = "{}"
$ws.Range("F5").Value = "{}"

$ws.Range("B6").Value = "{0: sala nr 11 | Dominik Kaczor | Informatyka}"
$ws.Range("C6").Value = "{0: sala nr 8 | Karolina Kamińska | Chemia}"
$ws.Range("D6").Value = "{0: sala nr 5 | Katarzyna Mazur | Fizyka}"
$ws.Range("E6").Value = "{}"
$ws.Range("F6").Value = "{}"

$ws.Range("B7").Value = "{0: sala nr 11 | Paweł Lewandowski | Matematyka}"
$ws.Range("C7").Value = "{0: sala nr 11 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("D7").Value = "{}"
$ws.Range("E7").Value = "{}"
$ws.Range("F7").Value = "{}"

$ws.Range("B8").Value = "{}"
$ws.Range("C8").Value = "{0: sala nr 6 | Piotr Wójcik | Biologia}"
$ws.Range("D8").Value = "{}"
$ws.Range("E8").Value = "{0: sala nr 3 | Lena Kowalska | Język angielski}"
$ws.Range("F8").Value = "{}"

$ws.Range("B9").Value = "{0: sala nr 11 | Jan Nowak | Język polski}"
$ws.Range("C9").Value = "{}"
$ws.Range("D9").Value = "{}"
$ws.Range("E9").Value = "{0: sala nr 7 | Karolina Kamińska | Chemia}"
$ws.Range("F9").Value = "{}"

$ws.Range("B10").Value = "{}"
$ws.Range("C10").Value = "{}"
$ws.Range("D10").Value = "{}"
$ws.Range("E10").Value = "{0: sala nr 1 | Dominik Kaczor | Informatyka}"
$ws.Range("F10").Value = "{}"

$ws.Range("B11").Value = "{}"
$ws.Range("C11").Value = "{}"
$ws.Range("D11").Value = "{0: sala nr 8 | Dominik Kaczor | Informatyka}"
$ws.Range("E11").Value = "{0: sala nr 9 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("F11").Value = "{0: sala nr 1 | Dominik Kaczor | Informatyka}"

$ws.Range("B12").Value = "{}"
$ws.Range("C12").Value = "{}"
$ws.Range("D12").Value = "{0: sala nr 2 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("E12").Value = "{0: sala nr 3 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("F12").Value = "{0: sala nr 7 | Jan Nowak | Język polski}"

$ws.Range("B13").Value = "{}"
$ws.Range("C13").Value = "{}"
$ws.Range("D13").Value = "{0: sala nr 3 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("E13").Value = "{0: sala nr 7 | Natalia Szymańska | Geografia}"
$ws.Range("F13").Value = "{0: sala nr 2 | Paweł Lewandowski | Matematyka}"

# Column width adjustments
$ws.Columns.Item(2).ColumnWidth = 47.8   # -> ~48.71 target, nearest achievable 48.67
$ws.Columns.Item(3).ColumnWidth = 55.8   # -> ~56.71 target, nearest achievable 56.67
$ws.Columns.Item(4).ColumnWidth = 50.8   # -> ~51.71 target, nearest achievable 51.67
$ws.Columns.Item(6).ColumnWidth = 46.8   # -> ~47.71 target, nearest achievable 47.67
